# Update Edn1-Ednra NATMI output with new TPM-derived expression values.
#
# The underlying single-cell TPM re-quantification changed the ligand
# (Edn1) average/total expression for the "ECs" sending cluster
# (columns G/H) and the receptor (Ednra) average/total expression for
# the "ECs" target cluster (columns M/N). Every other numeric column on
# the sheet (ligand/receptor derived-specificity columns I/J and O/P,
# the average/total edge expression weights Q/R, and the edge derived
# specificity columns S/T) is recomputed from those raw values:
#
#   I = J = G_row / sum(G over distinct sending clusters)
#   O = P = M_row / sum(M over distinct target clusters)
#   Q = G_row * M_row                       (average edge weight)
#   R = H_row * N_row                       (total edge weight)
#   S = Q_row / sum(Q over all rows)
#   T = R_row / sum(R over all rows)
#
# The values below are the results of that recalculation (matching the
# values produced by the updated NATMI scripts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 6.523120333333334
$ws.Range("H2").Value = 19.569361
$ws.Range("I2").Value = 0.7039464591847725
$ws.Range("J2").Value = 0.7039464591847725
$ws.Range("M2").Value = 2.067959333333333
$ws.Range("N2").Value = 6.203878
$ws.Range("O2").Value = 0.03751906553627715
$ws.Range("P2").Value = 0.03751906553627715
$ws.Range("Q2").Value = 13.48954757577311
$ws.Range("R2").Value = 121.405928181958
$ws.Range("S2").Value = 0.02641141333618373
$ws.Range("T2").Value = 0.02641141333618373

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 6.523120333333334
$ws.Range("H3").Value = 19.569361
$ws.Range("I3").Value = 0.7039464591847725
$ws.Range("J3").Value = 0.7039464591847725
$ws.Range("O3").Value = 0.1625861532004571
$ws.Range("P3").Value = 0.1625861532004571
$ws.Range("Q3").Value = 58.45597744535733
$ws.Range("R3").Value = 526.103797008216
$ws.Range("S3").Value = 0.1144519468579347
$ws.Range("T3").Value = 0.1144519468579347

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 6.523120333333334
$ws.Range("H4").Value = 19.569361
$ws.Range("I4").Value = 0.7039464591847725
$ws.Range("J4").Value = 0.7039464591847725
$ws.Range("O4").Value = 0.7998947812632657
$ws.Range("P4").Value = 0.7998947812632657
$ws.Range("Q4").Value = 287.5929491642161
$ws.Range("R4").Value = 2588.336542477945
$ws.Range("S4").Value = 0.563083098990654
$ws.Range("T4").Value = 0.563083098990654

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.2740496574363094
$ws.Range("J5").Value = 0.2740496574363095
$ws.Range("M5").Value = 2.067959333333333
$ws.Range("N5").Value = 6.203878
$ws.Range("O5").Value = 0.03751906553627715
$ws.Range("P5").Value = 0.03751906553627715
$ws.Range("Q5").Value = 5.251544125092444
$ws.Range("R5").Value = 47.263897125832
$ws.Range("S5").Value = 0.0102820870575472
$ws.Range("T5").Value = 0.0102820870575472

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.2740496574363094
$ws.Range("J6").Value = 0.2740496574363095
$ws.Range("O6").Value = 0.1625861532004571
$ws.Range("P6").Value = 0.1625861532004571
$ws.Range("S6").Value = 0.04455667958847258
$ws.Range("T6").Value = 0.04455667958847259

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.2740496574363094
$ws.Range("J7").Value = 0.2740496574363095
$ws.Range("O7").Value = 0.7998947812632657
$ws.Range("P7").Value = 0.7998947812632657
$ws.Range("S7").Value = 0.2192108907902896
$ws.Range("T7").Value = 0.2192108907902897

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.02200388337891808
$ws.Range("J8").Value = 0.02200388337891808
$ws.Range("M8").Value = 2.067959333333333
$ws.Range("N8").Value = 6.203878
$ws.Range("O8").Value = 0.03751906553627715
$ws.Range("P8").Value = 0.03751906553627715
$ws.Range("Q8").Value = 0.4216548401073333
$ws.Range("R8").Value = 3.794893560966
$ws.Range("S8").Value = 0.0008255651425462271
$ws.Range("T8").Value = 0.0008255651425462271

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.02200388337891808
$ws.Range("J9").Value = 0.02200388337891808
$ws.Range("O9").Value = 0.1625861532004571
$ws.Range("P9").Value = 0.1625861532004571
$ws.Range("S9").Value = 0.003577526754049767
$ws.Range("T9").Value = 0.003577526754049766

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.02200388337891808
$ws.Range("J10").Value = 0.02200388337891808
$ws.Range("O10").Value = 0.7998947812632657
$ws.Range("P10").Value = 0.7998947812632657
$ws.Range("S10").Value = 0.01760079148232209
$ws.Range("T10").Value = 0.01760079148232209
